$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.141.07'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '1.749.49'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.96'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5296'
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2812'
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06183'
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("D10").Value = '1.748.26'
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07176'
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.48'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6464'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.633'
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.48'
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9997'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9991'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '26.031.55'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.80'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006775'
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").Value = '1.971.01'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.338'
$ws.Range("E22").Value = '  +4.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.741'
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.238'
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.70'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.522'
$ws.Range("E26").Value = '  +1.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.30'
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.807'
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.99'
$ws.Range("E29").Value = '  +1.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08307'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.768'
$ws.Range("E31").Value = '  +2.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.643'
$ws.Range("E32").Value = '  +5.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04634'
$ws.Range("E33").Value = '  +4.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.642'
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.009'
$ws.Range("E35").Value = '  +2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6330'
$ws.Range("E36").Value = '  +3.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.706'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01625'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.981'
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9995'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.20'
$ws.Range("E41").Value = '  +1.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3934'
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7556'
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.077'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("E45").Value = '  +3.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.355'
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05343'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.72'
$ws.Range("E48").Value = '  +3.43%  '
$ws.Range("E49").Value = '  +3.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3480'
$ws.Range("E50").Value = '  +1.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.615'
$ws.Range("E51").Value = '  +0.84%  '
